# Swap the contents of columns B:AD between pairs of adjacent rows.
# Column A (the "id" index column) stays untouched in place.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(26, 27),
    @(32, 33),
    @(59, 60),
    @(65, 66),
    @(107, 108),
    @(115, 116),
    @(203, 204),
    @(225, 226),
    @(233, 234),
    @(248, 249),
    @(302, 303)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1`:AD$r1")
    $range2 = $ws.Range("B$r2`:AD$r2")

    $values1 = $range1.Value()
    $values2 = $range2.Value()

    $range1.Value = $values2
    $range2.Value = $values1
}
